$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "(419, 499, 117, 187)"
$ws.Range("B3").Value = "(634, 489, 87, 157)"
$ws.Range("B4").Value = "(912, 490, 89, 139)"
$ws.Range("B5").Value = "(1047, 480, 73, 132)"
$ws.Range("B6").Value = "(1196, 473, 47, 121)"
$ws.Range("B7").Value = "(1310, 476, 46, 101)"
